# New crime data collected
# Updates the 122nd Precinct weekly CompStat report:
#  - report volume/number and week-covering dates in the header
#  - the Crime Complaints table (rows 15-28) with the newly collected figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text (shared-string rich-text cells)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# ---------------------------------------------------------------------------
# Cells that switch from a numeric value to the text placeholders "0" / "***.*"
# (used by the report whenever a count is zero or a percentage is undefined).
# Format is copied from an existing placeholder cell so the cell keeps the
# same right-aligned style used throughout the table (style id 13).
# ---------------------------------------------------------------------------
function Set-PlaceholderText {
    param($cellRef, $text, $formatSourceRef)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

Set-PlaceholderText "C20" "0" "D23"
Set-PlaceholderText "C23" "0" "D23"
Set-PlaceholderText "G23" "0" "D23"
Set-PlaceholderText "H23" "***.*" "E23"
Set-PlaceholderText "D27" "0" "D23"
Set-PlaceholderText "E27" "***.*" "E23"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = -50
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -57.142857142857
$ws.Range("N15").Value = -70

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 20
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = -16.666666666666
$ws.Range("L16").Value = -20
$ws.Range("N16").Value = -82.608695652173

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 97
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 36.619718309859
$ws.Range("L17").Value = 22.784810126582
$ws.Range("M17").Value = 73.214285714285
$ws.Range("N17").Value = -11.009174311926

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 26.470588235294
$ws.Range("L18").Value = 10.256410256410
$ws.Range("M18").Value = -34.848484848484
$ws.Range("N18").Value = -92.081031307550

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -80
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -37.5
$ws.Range("I19").Value = 149
$ws.Range("J19").Value = 220
$ws.Range("K19").Value = -32.272727272727
$ws.Range("L19").Value = -24.747474747474
$ws.Range("M19").Value = -1.973684210526
$ws.Range("N19").Value = -55.522388059701

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -13.043478260869
$ws.Range("L20").Value = -48.717948717948
$ws.Range("M20").Value = -31.034482758620
$ws.Range("N20").Value = -98.183469573115

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -45
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = -8.571428571428
$ws.Range("I21").Value = 332
$ws.Range("J21").Value = 378
$ws.Range("K21").Value = -12.169312169312
$ws.Range("L21").Value = -13.541666666666
$ws.Range("M21").Value = -5.142857142857
$ws.Range("N21").Value = -84.997740623587

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("F23").Value = 9
$ws.Range("I23").Value = 28
$ws.Range("K23").Value = 366.666666666667
$ws.Range("L23").Value = 40
$ws.Range("M23").Value = 366.666666666667

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -53.571428571428
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -33.707865168539
$ws.Range("I24").Value = 403
$ws.Range("J24").Value = 421
$ws.Range("K24").Value = -4.275534441805
$ws.Range("L24").Value = -10.244988864142
$ws.Range("M24").Value = -37.904468412943

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -41.666666666666
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 3.030303030303
$ws.Range("I25").Value = 225
$ws.Range("J25").Value = 203
$ws.Range("K25").Value = 10.837438423645
$ws.Range("L25").Value = 15.979381443299

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 120
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 72
$ws.Range("I26").Value = 144
$ws.Range("J26").Value = 128
$ws.Range("K26").Value = 12.5
$ws.Range("L26").Value = 18.032786885245
$ws.Range("M26").Value = -37.117903930131

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = -66.666666666666
$ws.Range("L27").Value = -40

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 25
$ws.Range("L28").Value = 53.846153846153
